$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 309.12
$ws.Range("I15").Value = 309.12
$ws.Range("K15").Value = 927.36
$ws.Range("M15").Value = -758.36
$ws.Range("H44").Value = 3999.5
$ws.Range("J44").Value = 3999.5
$ws.Range("L44").Value = 3999.5
$ws.Range("N44").Value = -4923.5
$ws.Range("H112").Value = 3509.1853
$ws.Range("I112").Value = 899.75
$ws.Range("J112").Value = 3963
$ws.Range("K112").Value = 2699.25
$ws.Range("L112").Value = 11889
$ws.Range("M112").Value = -1591.25
$ws.Range("N112").Value = -14105
$ws.Range("H135").Value = 2591.4736
$ws.Range("I135").Value = 2319.2666
$ws.Range("J135").Value = 3612.25
$ws.Range("K135").Value = 20873.3994
$ws.Range("L135").Value = 32510.25
$ws.Range("M135").Value = -18338.3994
$ws.Range("N135").Value = -37580.25

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 3000
$ws.Range("J18").Value = 3000
$ws.Range("L18").Value = 3000
$ws.Range("N18").Value = -3644
$ws.Range("H74").Value = 1601.0667
$ws.Range("I74").Value = 1328.7273
$ws.Range("J74").Value = 2350
$ws.Range("K74").Value = 1328.7273
$ws.Range("L74").Value = 2350
$ws.Range("M74").Value = -454.7273
$ws.Range("N74").Value = -4098
$ws.Range("H77").Value = 1601.0667
$ws.Range("I77").Value = 1328.7273
$ws.Range("J77").Value = 2350
$ws.Range("K77").Value = 6643.636500000001
$ws.Range("L77").Value = 11750
$ws.Range("M77").Value = -2275.636500000001
$ws.Range("N77").Value = -20486
$ws.Range("H97").Value = 1578.3077
$ws.Range("I97").Value = 1247.9
$ws.Range("J97").Value = 2679.6667
$ws.Range("K97").Value = 1247.9
$ws.Range("L97").Value = 2679.6667
$ws.Range("M97").Value = -751.9000000000001
$ws.Range("N97").Value = -3671.6667
$ws.Range("H122").Value = 20975.906
$ws.Range("I122").Value = 22174.732
$ws.Range("J122").Value = 2993.5
$ws.Range("K122").Value = 66524.196
$ws.Range("L122").Value = 8980.5
$ws.Range("M122").Value = -64074.196
$ws.Range("N122").Value = -13880.5
$ws.Range("H132").Value = 2562.12
$ws.Range("I132").Value = 1574.7693
$ws.Range("J132").Value = 3631.75
$ws.Range("K132").Value = 4724.3079
$ws.Range("L132").Value = 10895.25
$ws.Range("M132").Value = -2194.3079
$ws.Range("N132").Value = -15955.25

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H134").Value = 19250
$ws.Range("I134").Value = 1300
$ws.Range("J134").Value = 37200
$ws.Range("K134").Value = 3900
$ws.Range("L134").Value = 111600
$ws.Range("M134").Value = -1365
$ws.Range("N134").Value = -116670

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 15000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 15000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 15000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -15348
$ws.Range("H31").Value = 2456.1714
$ws.Range("I31").Value = 2040.6957
$ws.Range("J31").Value = 3252.5
$ws.Range("K31").Value = 2040.6957
$ws.Range("L31").Value = 3252.5
$ws.Range("M31").Value = -1745.6957
$ws.Range("N31").Value = -3842.5
$ws.Range("H34").Value = 2456.1714
$ws.Range("I34").Value = 2040.6957
$ws.Range("J34").Value = 3252.5
$ws.Range("K34").Value = 2040.6957
$ws.Range("L34").Value = 3252.5
$ws.Range("M34").Value = -1838.6957
$ws.Range("N34").Value = -3656.5
$ws.Range("H50").Value = 9201
$ws.Range("J50").Value = 9201
$ws.Range("L50").Value = 9201
$ws.Range("N50").Value = -10451
$ws.Range("H68").Value = 17700.666
$ws.Range("J68").Value = 17700.666
$ws.Range("L68").Value = 17700.666
$ws.Range("N68").Value = -19198.666
$ws.Range("H71").Value = 17700.666
$ws.Range("J71").Value = 17700.666
$ws.Range("L71").Value = 53101.99800000001
$ws.Range("N71").Value = -60589.99800000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5102891
$ws.Range("J131").Value = 5556380
$ws.Range("L131").Value = 16669140
$ws.Range("N131").Value = -16679220

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7716.05
$ws.Range("I132").Value = 8354.588
$ws.Range("J132").Value = 4097.6665
$ws.Range("K132").Value = 25063.764
$ws.Range("L132").Value = 12292.9995
$ws.Range("M132").Value = -22533.764
$ws.Range("N132").Value = -17352.9995
$ws.Range("H135").Value = 47650
$ws.Range("J135").Value = 47650
$ws.Range("L135").Value = 47650
$ws.Range("N135").Value = -57790
$ws.Range("H136").Value = 9572.904
$ws.Range("J136").Value = 9572.904
$ws.Range("L136").Value = 28718.712
$ws.Range("N136").Value = -33818.712
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 57390
$ws.Range("J140").Value = 57390
$ws.Range("L140").Value = 57390
$ws.Range("N140").Value = -67750

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 300000
$ws.Range("I17").Value = 300000
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 300000
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -299830
$ws.Range("N17").ClearContents()
$ws.Range("H93").Value = 2674.2083
$ws.Range("I93").Value = 2363
$ws.Range("J93").Value = 2937.5386
$ws.Range("K93").Value = 2363
$ws.Range("L93").Value = 2937.5386
$ws.Range("M93").Value = -1115
$ws.Range("N93").Value = -5433.5386
$ws.Range("H136").Value = 2775.5334
$ws.Range("I136").Value = 3058
$ws.Range("J136").Value = 1998.75
$ws.Range("K136").Value = 9174
$ws.Range("L136").Value = 5996.25
$ws.Range("M136").Value = -6624
$ws.Range("N136").Value = -11096.25
$ws.Range("H137").Value = 45957.25
$ws.Range("J137").Value = 49609.668
$ws.Range("L137").Value = 49609.668
$ws.Range("N137").Value = -59809.668
$ws.Range("H139").Value = 79150
$ws.Range("J139").Value = 79150
$ws.Range("L139").Value = 79150
$ws.Range("N139").Value = -89430
$ws.Range("H140").Value = 87000
$ws.Range("J140").Value = 87000
$ws.Range("L140").Value = 87000
$ws.Range("N140").Value = -97360
$ws.Range("H141").Value = 59950
$ws.Range("J141").Value = 59950
$ws.Range("L141").Value = 59950
$ws.Range("N141").Value = -70310

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 22707.715
$ws.Range("J41").Value = 13159
$ws.Range("L41").Value = 13159
$ws.Range("N41").Value = -13939
$ws.Range("H132").Value = 38465510
$ws.Range("I132").Value = 45455748
$ws.Range("J132").Value = 19215.834
$ws.Range("K132").Value = 136367244
$ws.Range("L132").Value = 57647.50199999999
$ws.Range("M132").Value = -136364714
$ws.Range("N132").Value = -62707.50199999999
$ws.Range("H136").Value = 751.7059
$ws.Range("I136").Value = 648.5
$ws.Range("J136").Value = 1233.3334
$ws.Range("K136").Value = 1945.5
$ws.Range("L136").Value = 3700.0002
$ws.Range("M136").Value = 604.5
$ws.Range("N136").Value = -8800.0002
